$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 69

# New LeetCode problem row: Robot Bounded In Circle
$ws.Cells.Item($row, 1).Value = "Robot Bounded In Circle"
$ws.Cells.Item($row, 2).Value = "String"
$ws.Cells.Item($row, 3).Value = "No"
$ws.Cells.Item($row, 4).Value = "No"
$ws.Cells.Item($row, 5).Value = "Medium"
$ws.Cells.Item($row, 6).Value = "Easy"
$ws.Cells.Item($row, 7).Value = "1041 - Robot Bounded In Circle"

$ws.Hyperlinks.Add($ws.Cells.Item($row, 7), "1041%20-%20Robot%20Bounded%20In%20Circle")
$ws.Cells.Item($row, 7).Style = "Hyperlink"

# Extend the "Easy/Medium/Hard" cell-background conditional formatting down
# to the new row (originally D9:F68 -> D9:F69)
$cfRange = $ws.Range("D9:F69")
$cfHard = $cfRange.FormatConditions.Add(1, 3, '"Hard"')
$cfHard.Interior.Color = 255
$cfMedium = $cfRange.FormatConditions.Add(1, 3, '"Medium"')
$cfMedium.Interior.Color = 3243501
$cfEasy = $cfRange.FormatConditions.Add(1, 3, '"Easy"')
$cfEasy.Interior.Color = 4697456

# Extend the 4 list-based data validations down to the new row
$ws.Range("E2:F69").Validation.Delete()
$ws.Range("E2:F69").Validation.Add(3, 1, 1, '"Easy, Medium, Hard"')

$ws.Range("C2:C69").Validation.Delete()
$ws.Range("C2:C69").Validation.Add(3, 0, 1, '"Yes, No"')
$ws.Range("C2:C69").Validation.IgnoreBlank = $false

$ws.Range("D2:D69").Validation.Delete()
$ws.Range("D2:D69").Validation.Add(3, 1, 1, '"Yes, No"')

$ws.Range("B2:B69").Validation.Delete()
$ws.Range("B2:B69").Validation.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design"')

# Update selection / scroll position to match the author's final view
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L49").Select()
